$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.605.84'
$ws.Range('E2').Value = '  +1.12%  '
$ws.Range('D3').Value = '1.852.70'
$ws.Range('E3').Value = '  +1.18%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '263.83'
$ws.Range('E5').Value = '  +1.77%  '
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5262'
$ws.Range('E7').Value = '  +1.94%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3245'
$ws.Range('E8').Value = '  +0.81%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06804'
$ws.Range('E9').Value = '  +1.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.97'
$ws.Range('E10').Value = '  +1.90%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7843'
$ws.Range('E11').Value = '  +2.66%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07767'
$ws.Range('E12').Value = '  +1.20%  '
$ws.Range('D13').Value = '1.857.33'
$ws.Range('E13').Value = '  +1.24%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '88.78'
$ws.Range('E14').Value = '  +0.20%  '
$ws.Range('E15').Value = '  +0.65%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.000'
$ws.Range('E16').Value = '  -0.11%  '
$ws.Range('E17').Value = '  -0.24%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007961'
$ws.Range('E18').Value = '  +1.14%  '
$ws.Range('B19').Value = 'Dai'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.001'
$ws.Range('E19').Value = '  -0.05%  '
$ws.Range('D20').Value = '26.629.12'
$ws.Range('E20').Value = '  +1.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.641'
$ws.Range('E21').Value = '  +2.51%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.479'
$ws.Range('E22').Value = '  +0.81%  '
$ws.Range('E23').Value = '  +2.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '143.75'
$ws.Range('E24').Value = '  -0.72%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.174'
$ws.Range('E25').Value = '  -5.58%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.679'
$ws.Range('E26').Value = '  +1.84%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.01'
$ws.Range('E27').Value = '  +0.63%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '111.73'
$ws.Range('E28').Value = '  +1.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.191'
$ws.Range('E29').Value = '  +0.41%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08726'
$ws.Range('E30').Value = '  +0.25%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.04868'
$ws.Range('E32').Value = '  +0.99%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7237'
$ws.Range('E33').Value = '  +6.69%  '
$ws.Range('E34').Value = '  +0.92%  '
$ws.Range('E35').Value = '  +0.97%  '
$ws.Range('E36').Value = '  +0.72%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.273'
$ws.Range('E37').Value = '  +3.95%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01794'
$ws.Range('E38').Value = '  +0.91%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.4877'
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9018'
$ws.Range('E40').Value = '  +1.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '111.41'
$ws.Range('E41').Value = '  +0.27%  '
$ws.Range('E42').Value = '  -2.54%  '
$ws.Range('E43').Value = '  -0.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.688'
$ws.Range('E44').Value = '  +0.28%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4199'
$ws.Range('E45').Value = '  +0.50%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.05885'
$ws.Range('E46').Value = '  +0.13%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.000'
$ws.Range('E47').Value = '  -0.71%  '
$ws.Range('E48').Value = '  -1.04%  '
$ws.Range('E49').Value = '  -0.43%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.8909'
$ws.Range('E50').Value = '  +3.55%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '60.12'
$ws.Range('E51').Value = '  +1.98%  '
